# "Fixed the information matrix" -----------------------------------------
# Corrects several mis-entered 0/1 flags in the symptom/disease matrix,
# restyles a couple of header/data cells, adds a totals row (13) with
# conditional formatting that flags any 0 totals, and reserves some blank
# rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Correct mis-entered values in the matrix
# ---------------------------------------------------------------------
$ws.Range("AI3").Value = 1
$ws.Range("AI5").Value = 1
$ws.Range("AU6").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("AI9").Value = 1
$ws.Range("AI10").Value = 1

# ---------------------------------------------------------------------
# 2. Re-style a couple of header cells (AI1 / AU1) and clear the
#    wrap-text styling that had been applied down the AI / AU columns
# ---------------------------------------------------------------------
$ws.Range("AI1").Interior.ColorIndex = -4142
$ws.Range("AU1").Font.Name = "Calibri"

$ws.Range("AI7,AU7,AI8,AU8,AI9,AU9,AI10,AU10,AI11,AU11").WrapText = $false

# ---------------------------------------------------------------------
# 3. Totals row (row 13): SUM each column over rows 2-11
# ---------------------------------------------------------------------
$ws.Range("B13").Formula = "=SUM(B2:B11)"
$ws.Range("C13:AH13").Formula = "=SUM(C2:C11)"
$ws.Range("AI13").Formula = "=SUM(AI2:AI11)"
$ws.Range("AJ13:BK13").Formula = "=SUM(AJ2:AJ11)"

# Conditional formatting: flag any column total that is still 0
$totals = $ws.Range("B13:BK13")
$condition = $totals.FormatConditions.Add(1, 3, "0")
$condition.Font.Color = 393372
$condition.Interior.Color = 13551615

# ---------------------------------------------------------------------
# 4. Reserve some blank rows below the table
# ---------------------------------------------------------------------
$ws.Range("A19:A23").Value = ""
$ws.Range("A24:A28").Value = ""

# ---------------------------------------------------------------------
# 5. View tidy-up
# ---------------------------------------------------------------------
$ws.Range("B16").Select()
$excel.ActiveWindow.Zoom = 125

Write-Host "edit applied"
